$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "285656899 / Le Cong Huy" row (row 3) entirely, shifting the
# rows below it up. This removes the magic-number duplicate "nbk-vl" row
# and leaves the remaining data (334442222 / Nguyen Thi Tho / nbk-qn) as
# the new row 3.
$ws.Rows.Item(3).Delete()
